$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 209.18182
$ws.Range("I39").Value = 130.1
$ws.Range("K39").Value = 390.3
$ws.Range("M39").Value = -94.29999999999995
$ws.Range("H58").Value = 1998.3572
$ws.Range("I58").Value = 674.1667
$ws.Range("J58").Value = 2991.5
$ws.Range("K58").Value = 2022.5001
$ws.Range("L58").Value = 8974.5
$ws.Range("M58").Value = -1872.5001
$ws.Range("N58").Value = -9274.5
$ws.Range("H69").Value = 7567.1665
$ws.Range("H70").Value = 3799.4583
$ws.Range("J70").Value = 4906.4287
$ws.Range("L70").Value = 14719.2861
$ws.Range("N70").Value = -15259.2861
$ws.Range("H72").Value = 7567.1665
$ws.Range("H73").Value = 3799.4583
$ws.Range("J73").Value = 4906.4287
$ws.Range("L73").Value = 14719.2861
$ws.Range("N73").Value = -16591.2861
$ws.Range("H87").Value = 100001
$ws.Range("J87").Value = 100001
$ws.Range("L87").Value = 100001
$ws.Range("N87").Value = -102497
$ws.Range("H90").Value = 100001
$ws.Range("J90").Value = 100001
$ws.Range("L90").Value = 300003
$ws.Range("N90").Value = -312483
$ws.Range("H112").Value = 1754.6923
$ws.Range("J112").Value = 1367.8889
$ws.Range("L112").Value = 4103.6667
$ws.Range("N112").Value = -6319.6667
$ws.Range("H137").Value = 2036.1666
$ws.Range("I137").Value = 908.8333
$ws.Range("K137").Value = 2726.4999
$ws.Range("M137").Value = -176.4998999999998
$ws.Range("H138").Value = 2703.4375
$ws.Range("I138").Value = 881
$ws.Range("J138").Value = 4120.8887
$ws.Range("K138").Value = 2643
$ws.Range("L138").Value = 12362.6661
$ws.Range("M138").Value = 2497
$ws.Range("N138").Value = -22642.6661

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2022.8
$ws.Range("I2").Value = 2114.2222
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 2114.2222
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -2001.2222
$ws.Range("N2").Value = -1426
$ws.Range("H38").Value = 1000
$ws.Range("I38").Value = 1000
$ws.Range("K38").Value = 1000
$ws.Range("M38").Value = -533
$ws.Range("H45").Value = 2755.5715
$ws.Range("I45").Value = 2123.4
$ws.Range("K45").Value = 2123.4
$ws.Range("M45").Value = -1746.4
$ws.Range("H101").Value = 59500
$ws.Range("J101").Value = 59500
$ws.Range("L101").Value = 59500
$ws.Range("N101").Value = -65990
$ws.Range("H102").Value = 8338009.5
$ws.Range("I102").Value = 15626478
$ws.Range("J102").Value = 8331.286
$ws.Range("K102").Value = 15626478
$ws.Range("L102").Value = 8331.286
$ws.Range("M102").Value = -15624856
$ws.Range("N102").Value = -11575.286
$ws.Range("H116").Value = 2022.8
$ws.Range("I116").Value = 2114.2222
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 2114.2222
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = 179.7777999999998
$ws.Range("N116").Value = -5788

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2022.8
$ws.Range("I3").Value = 2114.2222
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 2114.2222
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = -2000.2222
$ws.Range("N3").Value = -1428
$ws.Range("H105").Value = 6994731
$ws.Range("I105").Value = 10102556
$ws.Range("J105").Value = 2124.75
$ws.Range("K105").Value = 10102556
$ws.Range("L105").Value = 2124.75
$ws.Range("M105").Value = -10100809
$ws.Range("N105").Value = -5618.75
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.92308
$ws.Range("J7").Value = 80.2
$ws.Range("L7").Value = 80.2
$ws.Range("N7").Value = -306.2
$ws.Range("H16").Value = 1950.6666
$ws.Range("I16").Value = 1819.5
$ws.Range("K16").Value = 1819.5
$ws.Range("M16").Value = -1532.5
$ws.Range("H22").Value = 1487.5
$ws.Range("I22").Value = 687.5714
$ws.Range("K22").Value = 687.5714
$ws.Range("M22").Value = -337.5714
$ws.Range("H99").Value = 2464
$ws.Range("I99").Value = 2459.6667
$ws.Range("K99").Value = 2459.6667
$ws.Range("M99").Value = -961.6667000000002
$ws.Range("H113").Value = 1950.6666
$ws.Range("I113").Value = 1819.5
$ws.Range("K113").Value = 1819.5
$ws.Range("M113").Value = 350.5
$ws.Range("H126").Value = 2464
$ws.Range("I126").Value = 2459.6667
$ws.Range("K126").Value = 7379.000100000001
$ws.Range("M126").Value = -4909.000100000001
$ws.Range("H132").Value = 2807.111
$ws.Range("I132").Value = 2807.111
$ws.Range("K132").Value = 8421.332999999999
$ws.Range("M132").Value = -5891.332999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1513.25
$ws.Range("I5").Value = 1632.6666
$ws.Range("J5").Value = 1393.8334
$ws.Range("K5").Value = 4897.9998
$ws.Range("L5").Value = 4181.5002
$ws.Range("M5").Value = -4785.9998
$ws.Range("N5").Value = -4405.5002
$ws.Range("H14").Value = 1401.3334
$ws.Range("I14").Value = 1401.3334
$ws.Range("K14").Value = 4204.0002
$ws.Range("M14").Value = -4031.0002
$ws.Range("H23").Value = 157.6
$ws.Range("J23").Value = 172
$ws.Range("L23").Value = 516
$ws.Range("N23").Value = -986
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H75").Value = 4901.5
$ws.Range("J75").Value = 5815.2
$ws.Range("L75").Value = 17445.6
$ws.Range("N75").Value = -19441.6
$ws.Range("H78").Value = 4901.5
$ws.Range("J78").Value = 5815.2
$ws.Range("L78").Value = 52336.8
$ws.Range("N78").Value = -62320.8
$ws.Range("H135").Value = 1513.25
$ws.Range("I135").Value = 1632.6666
$ws.Range("J135").Value = 1393.8334
$ws.Range("K135").Value = 14693.9994
$ws.Range("L135").Value = 12544.5006
$ws.Range("M135").Value = -12158.9994
$ws.Range("N135").Value = -17614.5006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 317.33334
$ws.Range("I2").Value = 147.3
$ws.Range("K2").Value = 147.3
$ws.Range("M2").Value = -34.30000000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1035.3077
$ws.Range("I22").Value = 997.5
$ws.Range("J22").Value = 1052.1111
$ws.Range("K22").Value = 997.5
$ws.Range("L22").Value = 1052.1111
$ws.Range("M22").Value = -702.5
$ws.Range("N22").Value = -1642.1111
$ws.Range("H27").Value = 1035.3077
$ws.Range("I27").Value = 997.5
$ws.Range("J27").Value = 1052.1111
$ws.Range("K27").Value = 997.5
$ws.Range("L27").Value = 1052.1111
$ws.Range("M27").Value = -890.5
$ws.Range("N27").Value = -1266.1111
$ws.Range("H40").Value = 1814.3611
$ws.Range("I40").Value = 2392.4546
$ws.Range("K40").Value = 2392.4546
$ws.Range("M40").Value = -2256.4546
$ws.Range("H55").Value = 1125.375
$ws.Range("I55").Value = 1717
$ws.Range("J55").Value = 856.4545000000001
$ws.Range("K55").Value = 1717
$ws.Range("L55").Value = 856.4545000000001
$ws.Range("M55").Value = -1544
$ws.Range("N55").Value = -1202.4545
$ws.Range("H82").Value = 4005.9092
$ws.Range("I82").Value = 683.3333
$ws.Range("J82").Value = 5251.875
$ws.Range("K82").Value = 683.3333
$ws.Range("L82").Value = 5251.875
$ws.Range("M82").Value = -322.3333
$ws.Range("N82").Value = -5973.875
$ws.Range("H85").Value = 4005.9092
$ws.Range("I85").Value = 683.3333
$ws.Range("J85").Value = 5251.875
$ws.Range("K85").Value = 683.3333
$ws.Range("L85").Value = 5251.875
$ws.Range("M85").Value = 564.6667
$ws.Range("N85").Value = -7747.875
$ws.Range("H93").Value = 2626.2856
$ws.Range("I93").Value = 2626.2856
$ws.Range("K93").Value = 2626.2856
$ws.Range("M93").Value = -1378.2856
$ws.Range("H132").Value = 1808
$ws.Range("I132").Value = 1808
$ws.Range("K132").Value = 5424
$ws.Range("M132").Value = -2894

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 88888
$ws.Range("J42").Value = 88888
$ws.Range("L42").Value = 88888
$ws.Range("N42").Value = -89644
$ws.Range("H62").Value = 8166.5
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 10999.75
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 10999.75
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -12247.75
$ws.Range("H65").Value = 8166.5
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 10999.75
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 54998.75
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -61238.75
$ws.Range("H107").Value = 66668170
$ws.Range("I107").Value = 111111620
$ws.Range("K107").Value = 333334860
$ws.Range("M107").Value = -333332940
$ws.Range("H113").Value = 1056.7142
$ws.Range("I113").Value = 849.125
$ws.Range("K113").Value = 2547.375
$ws.Range("M113").Value = -377.375
